$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) New row 31: Romania, appended after Slovakia (row 30), all zero values.
# ---------------------------------------------------------------------------
$ws.Range("A31").Value = "Romania"
$ws.Range("B31:E31").Value = 0
$ws.Range("F31:J31").Value = 0

# ---------------------------------------------------------------------------
# 2) Apply the existing "red" cell style (the one already used on cells such
#    as A3 / F11:J11) to the newly-touched ranges. Copy format from a cell
#    that already carries that style so Excel re-uses the existing style
#    index instead of minting a new one.
# ---------------------------------------------------------------------------
$ws.Range("A3").Copy() | Out-Null

$ws.Range("F3:J3").PasteSpecial(-4122)
$ws.Range("F4:J4").PasteSpecial(-4122)
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("F12:J12").PasteSpecial(-4122)
$ws.Range("F14:J14").PasteSpecial(-4122)
$ws.Range("A15").PasteSpecial(-4122)
$ws.Range("F17:J17").PasteSpecial(-4122)
$ws.Range("F20:J20").PasteSpecial(-4122)
$ws.Range("F21:J21").PasteSpecial(-4122)
$ws.Range("F27:J27").PasteSpecial(-4122)
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("F29:J29").PasteSpecial(-4122)
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("F30:J30").PasteSpecial(-4122)
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("F31:J31").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3) Row 28 (Serbia) switches to the other existing custom style (the one
#    already used on A26) across the whole data range.
# ---------------------------------------------------------------------------
$ws.Range("A26").Copy() | Out-Null
$ws.Range("A28:J28").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Row 24 (Finland) loses its red style on A24 - revert back to the plain
#    default style by copying format from a default-styled cell.
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy() | Out-Null
$ws.Range("A24").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 5) Corrected policy values: several countries flip from 0 to 1 for
#    years 2018-2022 (columns F:J).
# ---------------------------------------------------------------------------
$ws.Range("F11:J11").Value = 1
$ws.Range("F24:J24").Value = 1
$ws.Range("F28:J28").Value = 1

# ---------------------------------------------------------------------------
# 6) Sheet view bookkeeping - matches the refreshed selection/scroll state.
# ---------------------------------------------------------------------------
$ws.Range("M24").Select()
